$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.386.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.118.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.15%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -8.18%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.116.08"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.86%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.627.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.499.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.118.35"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.14"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.71"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -7.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.92%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.117"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -11.27%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -12.02%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.66"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.10"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.95%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.27%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.69"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0751"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "452.26"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -13.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0390"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.05%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.836.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.90%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.262"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -9.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.99"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.17%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.49%  "
